$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Socio "ferraro / Franco Valentin" (row 2) settled their balance and
# their membership was renewed with a new due date:
#  - Haber (H2) goes from 4500 to 0 and is highlighted green (cleared),
#    matching the green "paid/regular" highlight already used elsewhere
#    in the sheet (e.g. H3:H5, H7, H9).
#  - Vencimiento_Formateada (J2) moves from 22/02/2024 to 25/01/2024.
$ws.Range("H2").Value = "0"
$ws.Range("H2").Interior.Color = 582476
$ws.Range("J2").Value = "25/01/2024"
